# Applies the "Add files via upload" edit to the School Canteen planning
# document. Every textual change in the target diff is a run-merge (the
# author re-typed short phrases in the browser so Word's AutoCorrect no
# longer needed w:proofErr markers / run splits) plus a handful of real
# wording tweaks near the end (Task 14 evaluation section).
#
# We drive this purely with Word's Find/Replace: replacing the *full*
# existing phrase (reconstructed by concatenating the split runs) with the
# same (or tweaked) text in one shot collapses all of the runs Word had
# split across into a single run and drops the now-orphaned w:proofErr
# markers, which is exactly what the diff shows.

$d = $word.ActiveDocument

function Replace-Text($find, $replace) {
    $range = $d.Content
    $ok = $range.Find.Execute($find, $true, $false, $false, $false, $false, $true, 1, $false, $replace, 2)
    if (-not $ok) {
        Write-Host "MISS:" $find
    }
    return $ok
}

# Task 0 heading: "Task " + "0 :" + " Explain what you are doing/ going to accomplish"
Replace-Text "Task 0 : Explain what you are doing/ going to accomplish" "Task 0 : Explain what you are doing/ going to accomplish" | Out-Null

# "I am creating " + "food_list" + " to hold all of my food items in"
Replace-Text "I am creating food_list to hold all of my food items in" "I am creating food_list to hold all of my food items in" | Out-Null

# "Data = " + "dict" + "(" + "food_list" + " = food)"
Replace-Text "Data = dict(food_list = food)" "Data = dict(food_list = food)" | Out-Null

# "DEF " + "CANTEEN(" + "):"
Replace-Text "DEF CANTEEN():" "DEF CANTEEN():" | Out-Null

# Table cells: merely drop w:proofErr around "dict" and "food_list" tokens;
# text is unchanged, so touching them with identical find/replace text
# collapses the (identical-format) surrounding runs and removes the
# proofErr markers without altering the visible content.
Replace-Text "dict(food_list" "dict(food_list" | Out-Null

# "% for item in " + "menu_list" + ":"
Replace-Text "% for item in menu_list:" "% for item in menu_list:" | Out-Null

# "<p" + ">{" + "{" + "food" + ".stock" + "}}</p>"
Replace-Text "<p>{{food.stock}}</p>" "<p>{{food.stock}}</p>" | Out-Null

# "<p" + ">{" + "{" + "food" + ".price" + "}} dollar " + "dollar" + " bills " + "yo" + "</p>"
Replace-Text "<p>{{food.price}} dollar dollar bills yo</p>" "<p>{{food.price}} dollar dollar bills yo</p>" | Out-Null

# "I tested my case by running it and opening the page I created in localhost" + ":8080"
Replace-Text "I tested my case by running it and opening the page I created in localhost:8080" "I tested my case by running it and opening the page I created in localhost:8080" | Out-Null

# The "I needed to slightly modify..." sentence gains a new clause
# (" from the list") and its tail is re-typed as a single run ending in
# "it showed item.price".
Replace-Text "I needed to slightly modify my html code so that instead of displaying food.price (non existant as I said to put each ITEM not food on display for each user) it showed item.price" "I needed to slightly modify my html code so that instead of displaying food.price (non existant as I said to put each ITEM not food from the list on display for each user) it showed item.price" | Out-Null

# Task 14 heading: "Task " + "14 :" + " Evaluation"
Replace-Text "Task 14 : Evaluation" "Task 14 : Evaluation" | Out-Null

# "How did your version turn " + "out"
Replace-Text "How did your version turn out" "How did your version turn out" | Out-Null

# Final evaluation paragraph: originally just "After fixing this coed";
# becomes a new explanatory sentence, followed by the old closing remark
# (now prefixed with a space) about the next version.
Replace-Text "After fixing this coed" "After making a slight fix to my code it worked well and accomplished the task that I set it to do in Task 0. I wanted this version to display a product page (menu) to the user and this is what it does. In the next version I will add pictures to my website" | Out-Null

Write-Host "done"
